$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$excel.ActiveWindow.ScrollRow = 71
$excel.ActiveWindow.ScrollColumn = 2
